$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column: values are numeric-looking strings stored as text in the
# source workbook (t="inlineStr"). A leading apostrophe forces Excel to
# keep them as text (General format) instead of auto-converting to numbers,
# matching the original cell formatting/type.
# E column: values already contain spaces/"%" so they remain text naturally.

$ws.Range("D2").Value = "'62.913.35"
$ws.Range("E2").Value = '  -0.37%  '
$ws.Range("D3").Value = "'2.546.13"
$ws.Range("E3").Value = '  +3.20%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'567.48"
$ws.Range("E5").Value = '  +0.02%  '
$ws.Range("D6").Value = "'147.07"
$ws.Range("E6").Value = '  +2.38%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").Value = "'2.544.33"
$ws.Range("E9").Value = '  +3.20%  '
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("E11").Value = '  -2.37%  '
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = "'27.28"
$ws.Range("E14").Value = '  +2.77%  '
$ws.Range("D15").Value = "'3.000.61"
$ws.Range("E15").Value = '  +3.23%  '
$ws.Range("D16").Value = "'62.885.93"
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = "'2.551.51"
$ws.Range("E18").Value = '  +3.48%  '
$ws.Range("E19").Value = '  +2.10%  '
$ws.Range("D20").Value = "'336.83"
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("D22").Value = "'6.74"
$ws.Range("E22").Value = '  -1.01%  '
$ws.Range("E23").Value = '  +0.10%  '
$ws.Range("D24").Value = "'65.25"
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("E25").Value = '  +8.71%  '
$ws.Range("E26").Value = '  -3.18%  '
$ws.Range("D27").Value = "'1.50"
$ws.Range("E27").Value = '  +11.47%  '
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("E29").Value = '  +2.79%  '
$ws.Range("D30").Value = "'7.27"
$ws.Range("E30").Value = '  +6.43%  '
$ws.Range("D31").Value = "'0.0₃0813"
$ws.Range("E31").Value = '  +0.23%  '
$ws.Range("E32").Value = '  -0.22%  '
$ws.Range("D33").Value = "'178.17"
$ws.Range("E33").Value = '  +0.82%  '
$ws.Range("E34").Value = '  +2.75%  '
$ws.Range("D35").Value = "'406.89"
$ws.Range("E35").Value = '  +10.17%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").Value = "'18.96"
$ws.Range("E37").Value = '  +0.12%  '
$ws.Range("E38").Value = '  -0.02%  '
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("D40").Value = "'1.77"
$ws.Range("E40").Value = '  +3.62%  '
$ws.Range("D41").Value = "'0.999"
$ws.Range("E41").Value = '  -0.05%  '
$ws.Range("D42").Value = "'39.11"
$ws.Range("E42").Value = '  -3.22%  '
$ws.Range("D43").Value = "'153.08"
$ws.Range("E43").Value = '  +1.38%  '
$ws.Range("E44").Value = '  +0.41%  '
$ws.Range("D45").Value = "'20.77"
$ws.Range("E45").Value = '  +0.81%  '
$ws.Range("D46").Value = "'0.607"
$ws.Range("E46").Value = '  +1.56%  '
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").Value = "'0.0518"
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").Value = "'0.0236"
$ws.Range("E49").Value = '  +4.25%  '
$ws.Range("D50").Value = "'18.23"
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("E51").Value = '  +0.54%  '
